$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new row at position 68 - this pushes the existing rows 68-106
# down to 69-107, and expands the sheet dimension from A1:R106 to A1:R107.
$ws.Rows.Item(68).Insert()

# Populate the newly-inserted row 68 with the new weekly "Ajo" price record.
$ws.Range("A68").Value = 11
$ws.Range("B68").Value = 'Vega Monumental Concepción'
$ws.Range("C68").Value = 'Bíobío'
$ws.Range("D68").Value = 44518
$ws.Range("E68").Value = 8
$ws.Range("F68").Value = 100112003
$ws.Range("G68").Value = 'Ajo'
$ws.Range("H68").Value = 'Chino'
$ws.Range("I68").Value = 'Primera'
$ws.Range("J68").Value = 310
$ws.Range("K68").Value = 17000
$ws.Range("L68").Value = 18000
$ws.Range("M68").Value = 17516
$ws.Range("N68").Value = '$/caja 10 kilos'
$ws.Range("O68").Value = 'China'
$ws.Range("P68").Value = 1752
$ws.Range("Q68").Value = 10
$ws.Range("R68").Value = 'Hortaliza'
